# API: new api for updating product price
# Rename the sheet, update headers, and append a new data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "Sheet" to "Master"
$ws.Name = "Master"

# --- Update header row (row 1) ---
$ws.Range("A1").Value = "Barcode"
$ws.Range("B1").Value = "Product name"
$ws.Range("C1").Value = "Supplier Code"
$ws.Range("D1").Value = "Supplier Name"
$ws.Range("E1").Value = "Quantity"
# New header cell E1 needs the same bold header styling as the rest of row 1
$ws.Range("E1").Font.Bold = $true

# --- Add new data row (row 2) ---
# Force text storage first so numeric-looking values (barcode, quantity)
# are written as text instead of being auto-converted to numbers, then
# reset the style back to Normal so no lingering number format remains
# on the cells (matches a plain, unstyled data row).
$ws.Range("A2:E2").NumberFormat = "@"
$ws.Range("A2").Value = "770795005596"
$ws.Range("B2").Value = "DAN-D-PAK HẠT ĐIỀU KHÔNG MUỐI 50G/1 GÓI"
$ws.Range("C2").Value = "C0017"
$ws.Range("D2").Value = "CÔNG TY CỔ PHẦN THỰC PHẨM DÂN ÔN"
$ws.Range("E2").Value = "0"
$ws.Range("A2:E2").Style = "Normal"
